# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The worker "AUGUSTO RAFAEL MELO PAJARO" statement periods are re-sorted
# descending (2210 -> 1908), and the worker "ALEX DE JESUS NAVARRO MARTINEZ"
# rows are moved to the top of the table (sorted Name asc, Period desc),
# with the F16/F18/F56 "Valor Mora" figures updated to match the database.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=16; C="73556042";  D="ALEX DE JESUS NAVARRO MARTINEZ"; E="1909"; F=50666;  G=1900000},
    @{Row=17; C="73556042";  D="ALEX DE JESUS NAVARRO MARTINEZ"; E="1908"; F=25333;  G=1900000},
    @{Row=18; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2210"; F=90000;   G=2500000},
    @{Row=19; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2209"; F=100000;  G=2500000},
    @{Row=20; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2208"; F=100000;  G=2500000},
    @{Row=21; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2207"; F=100000;  G=2500000},
    @{Row=22; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2206"; F=100000;  G=2500000},
    @{Row=23; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2205"; F=100000;  G=2500000},
    @{Row=24; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2204"; F=100000;  G=2500000},
    @{Row=25; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2203"; F=100000;  G=2500000},
    @{Row=26; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2202"; F=100000;  G=2500000},
    @{Row=27; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2201"; F=100000;  G=2500000},
    @{Row=28; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2112"; F=100000;  G=2500000},
    @{Row=29; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2111"; F=100000;  G=2500000},
    @{Row=30; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2110"; F=100000;  G=2500000},
    @{Row=31; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2109"; F=100000;  G=2500000},
    @{Row=32; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2108"; F=100000;  G=2500000},
    @{Row=33; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2107"; F=100000;  G=2500000},
    @{Row=34; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2106"; F=100000;  G=2500000},
    @{Row=35; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2105"; F=100000;  G=2500000},
    @{Row=36; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2104"; F=100000;  G=2500000},
    @{Row=37; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2103"; F=100000;  G=2500000},
    @{Row=38; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2102"; F=100000;  G=2500000},
    @{Row=39; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2101"; F=100000;  G=2500000},
    @{Row=40; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2012"; F=100000;  G=2500000},
    @{Row=41; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2011"; F=100000;  G=2500000},
    @{Row=42; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2010"; F=100000;  G=2500000},
    @{Row=43; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2009"; F=100000;  G=2500000},
    @{Row=44; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2008"; F=100000;  G=2500000},
    @{Row=45; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2007"; F=100000;  G=2500000},
    @{Row=46; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2006"; F=100000;  G=2500000},
    @{Row=47; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2005"; F=100000;  G=2500000},
    @{Row=48; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2004"; F=100000;  G=2500000},
    @{Row=49; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2003"; F=100000;  G=2500000},
    @{Row=50; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2002"; F=100000;  G=2500000},
    @{Row=51; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="2001"; F=100000;  G=2500000},
    @{Row=52; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="1912"; F=100000;  G=2500000},
    @{Row=53; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="1911"; F=100000;  G=2500000},
    @{Row=54; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="1910"; F=100000;  G=2500000},
    @{Row=55; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="1909"; F=100000;  G=2500000},
    @{Row=56; C="73134285";  D="AUGUSTO RAFAEL MELO PAJARO";      E="1908"; F=33333;   G=2500000}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value2 = "CC"
    $ws.Cells.Item($r, 3).Value2 = $item.C
    $ws.Cells.Item($r, 4).Value2 = $item.D
    $ws.Cells.Item($r, 5).Value2 = $item.E
    $ws.Cells.Item($r, 6).Value2 = $item.F
    $ws.Cells.Item($r, 7).Value2 = $item.G
}
